# Smoke / MaestroData.xlsx update
# - Refresh the smoke-test seed data on sheet "DatosCuenta":
#     * rename test values "PruebaCuenta" -> "PruebaSmoke" and "Apellido" -> "ApellidoSmoke"
#     * bump the document/street numbers (27100101 -> 27100102, 106 -> 107)
# - Move the active sheet/selection from "DatosAP" back to "DatosCuenta" (D3)

$wb = $excel.ActiveWorkbook

$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "PruebaSmoke"
$wsCuenta.Range("B2").Value = "ApellidoSmoke"
$wsCuenta.Range("C2").Value = 27100102
$wsCuenta.Range("D2").Value = 107

# Make DatosCuenta the active sheet/tab again (it was DatosAP before) and
# move the selection to D3, matching the refreshed dataset.
$wsCuenta.Activate()
$wsCuenta.Range("D3").Select()
